# Update the 丽水-漫展信息 workbook to the newer scraped snapshot.
#
# The "展览" (exhibitions) sheet and its mirror "全部类型" (all-types) sheet
# both listed 8 rows (1 header + 7 events) starting with the now-stale
# "丽水·龙泉ACG动漫游戏博览会" (2024-07-20) entry. That entry has dropped out
# of the feed, so every later event shifts up one row (the last event,
# "丽水·LZ栗子动漫游戏嘉年华" on 2024-09-16, becomes the new final row with
# nothing appended after it), and three "想去人数" (interest count) figures
# ticked up slightly from a later crawl of the same event pages.
#
# Net effect per affected sheet: delete row 2, renumber the leading index
# column, and bump the three interest counts that changed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the stale first event row; everything below shifts up one row.
    $ws.Rows(2).Delete()

    # Column A holds a 0-based running index (header row = 0, first event
    # row = 1, ...). After the shift it needs to be renumbered to stay
    # sequential.
    for ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the "想去人数" (interest count) figures that moved and were
    # also updated by the later crawl.
    $ws.Cells.Item(3, 6).Value = 490   # 丽水·第四届HP国风动漫游戏嘉年华
    $ws.Cells.Item(4, 6).Value = 162   # 丽水·樱卡动漫游戏嘉年华
    $ws.Cells.Item(6, 6).Value = 674   # 丽水·AEO纯白礼赞动漫嘉年华
}
